# "Change import 1 to multiple files"
# Rename the header cell A1 ("name" -> "table_name") on both sheets, and
# switch the active sheet/selection from "list_files_to_import_db" back to
# "list_folders_to_import_db" (first tab, cell A2 selected on both sheets).

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("list_folders_to_import_db")
$ws2 = $wb.Worksheets.Item("list_files_to_import_db")

# Rename the "name" column header to "table_name" on both sheets.
$ws1.Range("A1").Value = "table_name"
$ws2.Range("A1").Value = "table_name"

# Update sheet2's selection (it stays the inactive sheet) to A2.
$ws2.Activate()
$ws2.Range("A2").Select()

# Make sheet1 the active sheet, with A2 selected.
$ws1.Activate()
$ws1.Range("A2").Select()
